# Remove the duplicate trailing "Laptop Environment for Lab" slide
# (a leftover copy of slide 3) together with its notes page.
$p = $ppt.ActivePresentation
$lastIndex = $p.Slides.Count
$s = $p.Slides.Item($lastIndex)
$s.Delete()
